$wb = $excel.ActiveWorkbook

# --- Update the Metrics sheet values (B2:B13) ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 447880.55
$metrics.Range("B3").Value  = 362291.19
$metrics.Range("B4").Value  = 141236.46
$metrics.Range("B5").Value  = 17841
$metrics.Range("B6").Value  = 4367131.4700000007
$metrics.Range("B7").Value  = 3689818.6700000004
$metrics.Range("B8").Value  = 1270602.1399999999
$metrics.Range("B9").Value  = 169001
$metrics.Range("B10").Value = 32832455.260000002
$metrics.Range("B11").Value = 30965040.229999997
$metrics.Range("B12").Value = 11552311.050000001
$metrics.Range("B13").Value = 1266628

# Formulas on the "today" sheet reference these Metrics cells (and chain
# through E/F columns), so they recalc automatically once the workbook
# recalculates after this script finishes.

# --- Update selection on Metrics sheet ---
$metrics.Activate() | Out-Null
$metrics.Range("E8").Select() | Out-Null

# --- Update selection on the "today" sheet and make it the active tab ---
$today = $wb.Worksheets.Item("today")
$today.Activate() | Out-Null
$today.Range("F9").Select() | Out-Null
